$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet1: only the selected cell changes ----
$ws1.Range("I11").Select()

# ---- Sheet2 ("Summary"): insert two rows at top for a new title row ----
$ws2.Rows.Item(1).Insert()
$ws2.Rows.Item(1).Insert()

# New title in A1, using the same red-bold "title" style used elsewhere in the workbook (s=3)
$ws1.Range("A2").Copy()
$ws2.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A1").Value = "Using LR, one hot encoding and ngram(1,2)"

# Fill in the previously-missing Accuracy(%) values (column D)
$ws2.Range("D5").Value = 94.8
$ws2.Range("D6").Value = 86.7
$ws2.Range("D7").Value = 89.7
$ws2.Range("D8").Value = 94
$ws2.Range("D9").Value = 93.6
$ws2.Range("D10").Value = 89.7
$ws2.Range("D11").Value = 93.7

$ws2.Range("D13").Value = 95
$ws2.Range("D14").Value = 94
$ws2.Range("D15").Value = 94
$ws2.Range("D16").Value = 93.7
$ws2.Range("D17").Value = 94.1
$ws2.Range("D18").Value = 93.8
$ws2.Range("D19").Value = 94.1

# Widen column A to fit the new title, drop the old bestFit autosize flag
$ws2.Columns.Item(1).ColumnWidth = 38.6667

# Update the saved selection for this sheet
$ws2.Range("D13").Select()
